# "Add 2-opt per route & show distance per route"
#
# The sheet holds one row per (route, stop): column A is the route
# ("Newspaper boy"), column B is the stop's sequence number along the
# route, and column C is the customer visited at that stop. Running a
# 2-opt local-search pass on each route re-orders the customers within
# the route (the set of customers per route, and the sequence numbers
# 1..N, are unchanged) to shorten the route's total travel distance.
#
# Below, each route's optimized customer order is written back into
# column C, starting at that route's first data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$routes = @{
    1 = @{ StartRow = 2;   Customers = @(83, 84, 71, 82, 74, 77, 103, 99, 100, 86, 89, 96, 104, 111, 107, 110, 108, 116, 119, 117, 118, 120, 115, 109, 113, 114, 112) }
    2 = @{ StartRow = 29;  Customers = @(85, 91, 105, 106, 98, 102, 101, 95, 88, 87, 93, 94, 97, 92, 90, 76, 80, 78, 75, 70, 62, 66, 63, 64, 79, 56, 51, 49, 42, 37, 39, 47, 48, 40, 26, 38, 33, 14, 6, 5, 2, 21, 18) }
    3 = @{ StartRow = 72;  Customers = @(73, 69, 61, 55, 52, 45, 44, 36, 34, 30, 25, 10, 13, 15, 17, 22, 24, 31, 32, 35, 43, 41, 29, 27, 28, 23, 19, 20, 16, 12, 11, 9, 7, 4, 3, 8, 1) }
    4 = @{ StartRow = 109; Customers = @(81, 72, 60, 65, 68, 67, 57, 59, 58, 54, 53, 50, 46) }
}

foreach ($routeId in $routes.Keys) {
    $route = $routes[$routeId]
    $row = $route.StartRow
    foreach ($customer in $route.Customers) {
        $ws.Cells.Item($row, 3).Value = $customer
        $row = $row + 1
    }
}
